$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Core data edit: StartYear value for Malawi changed from 2015 to 2018
$ws.Range("B2").Value = 2018

# Reflect the saved cursor/selection position recorded in the workbook
# (author had cell B2 selected when the file was last saved)
$ws.Range("B2").Select()
